$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regional/legal-category breakdown rows added for the 2020-09-22 data refresh.
# Columns C (nombre_aides) and D (montant_total) are stored as text in this sheet,
# so we force a Text number format while writing then restore Normal style afterward
# to avoid leaving a stray style index on the cell.
$rng = $ws.Range("C2:D2")
$rng.NumberFormat = "@"
$ws.Range("C2").Value = "1248"
$ws.Range("D2").Value = "3007291.70"
$rng.Style = "Normal"

$rng = $ws.Range("C4:D4")
$rng.NumberFormat = "@"
$ws.Range("C4").Value = "1100"
$ws.Range("D4").Value = "4372820.27"
$rng.Style = "Normal"

$rng = $ws.Range("C6:D6")
$rng.NumberFormat = "@"
$ws.Range("C6").Value = "737"
$ws.Range("D6").Value = "2610308.52"
$rng.Style = "Normal"

$rng = $ws.Range("C14:D14")
$rng.NumberFormat = "@"
$ws.Range("C14").Value = "243"
$ws.Range("D14").Value = "647219.00"
$rng.Style = "Normal"

$rng = $ws.Range("C16:D16")
$rng.NumberFormat = "@"
$ws.Range("C16").Value = "522"
$ws.Range("D16").Value = "1986149.55"
$rng.Style = "Normal"

$rng = $ws.Range("C17:D17")
$rng.NumberFormat = "@"
$ws.Range("C17").Value = "149"
$ws.Range("D17").Value = "472601.67"
$rng.Style = "Normal"

$rng = $ws.Range("C22:D22")
$rng.NumberFormat = "@"
$ws.Range("C22").Value = "354"
$ws.Range("D22").Value = "1468066.26"
$rng.Style = "Normal"

$rng = $ws.Range("C26:D26")
$rng.NumberFormat = "@"
$ws.Range("C26").Value = "130"
$ws.Range("D26").Value = "313001.13"
$rng.Style = "Normal"

$rng = $ws.Range("C27:D27")
$rng.NumberFormat = "@"
$ws.Range("C27").Value = "142"
$ws.Range("D27").Value = "498516.96"
$rng.Style = "Normal"

$rng = $ws.Range("C28:D28")
$rng.NumberFormat = "@"
$ws.Range("C28").Value = "121"
$ws.Range("D28").Value = "379045.60"
$rng.Style = "Normal"

$rng = $ws.Range("C30:D30")
$rng.NumberFormat = "@"
$ws.Range("C30").Value = "324"
$ws.Range("D30").Value = "844556.89"
$rng.Style = "Normal"

$rng = $ws.Range("C32:D32")
$rng.NumberFormat = "@"
$ws.Range("C32").Value = "639"
$ws.Range("D32").Value = "2946876.47"
$rng.Style = "Normal"

$rng = $ws.Range("C34:D34")
$rng.NumberFormat = "@"
$ws.Range("C34").Value = "424"
$ws.Range("D34").Value = "1540899.35"
$rng.Style = "Normal"

$rng = $ws.Range("C47:D47")
$rng.NumberFormat = "@"
$ws.Range("C47").Value = "507"
$ws.Range("D47").Value = "1428838.34"
$rng.Style = "Normal"

$rng = $ws.Range("C49:D49")
$rng.NumberFormat = "@"
$ws.Range("C49").Value = "739"
$ws.Range("D49").Value = "3414717.57"
$rng.Style = "Normal"

$rng = $ws.Range("C50:D50")
$rng.NumberFormat = "@"
$ws.Range("C50").Value = "516"
$ws.Range("D50").Value = "2083483.62"
$rng.Style = "Normal"

$rng = $ws.Range("C80:D80")
$rng.NumberFormat = "@"
$ws.Range("C80").Value = "979"
$ws.Range("D80").Value = "3731818.58"
$rng.Style = "Normal"

$rng = $ws.Range("C81:D81")
$rng.NumberFormat = "@"
$ws.Range("C81").Value = "547"
$ws.Range("D81").Value = "1972881.43"
$rng.Style = "Normal"

$rng = $ws.Range("C83:D83")
$rng.NumberFormat = "@"
$ws.Range("C83").Value = "37"
$ws.Range("D83").Value = "166180.27"
$rng.Style = "Normal"

$rng = $ws.Range("C84:D84")
$rng.NumberFormat = "@"
$ws.Range("C84").Value = "667"
$ws.Range("D84").Value = "1684707.33"
$rng.Style = "Normal"

$rng = $ws.Range("C87:D87")
$rng.NumberFormat = "@"
$ws.Range("C87").Value = "1012"
$ws.Range("D87").Value = "3812141.44"
$rng.Style = "Normal"

$rng = $ws.Range("C88:D88")
$rng.NumberFormat = "@"
$ws.Range("C88").Value = "724"
$ws.Range("D88").Value = "2319622.16"
$rng.Style = "Normal"

$rng = $ws.Range("C90:D90")
$rng.NumberFormat = "@"
$ws.Range("C90").Value = "34"
$ws.Range("D90").Value = "121571.23"
$rng.Style = "Normal"

$rng = $ws.Range("C97:D97")
$rng.NumberFormat = "@"
$ws.Range("C97").Value = "885"
$ws.Range("D97").Value = "2152465.94"
$rng.Style = "Normal"

$rng = $ws.Range("C100:D100")
$rng.NumberFormat = "@"
$ws.Range("C100").Value = "1233"
$ws.Range("D100").Value = "4447359.41"
$rng.Style = "Normal"

$rng = $ws.Range("C102:D102")
$rng.NumberFormat = "@"
$ws.Range("C102").Value = "1170"
$ws.Range("D102").Value = "3822979.27"
$rng.Style = "Normal"

$rng = $ws.Range("C104:D104")
$rng.NumberFormat = "@"
$ws.Range("C104").Value = "66"
$ws.Range("D104").Value = "256647.45"
$rng.Style = "Normal"
